$wb = $excel.ActiveWorkbook

# Sheets affected: "展览" (first tab) and "全部类型" (fourth tab)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 2966
    $ws.Range("F5").Value = 155
    $ws.Range("F14").Value = 336
    $ws.Range("F20").Value = 3043
}

$wb.Save()
